$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" header in H1 — same value/style as the rest of the header row (B1:G1)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Label values: 0 for Control patients (rows 2-6, 12-16), 1 for MDD patients (rows 7-11, 17-21)
$labels = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0
    7 = 1; 8 = 1; 9 = 1; 10 = 1; 11 = 1
    12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0
    17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}
